$wb = $excel.ActiveWorkbook

$loginSheet = $wb.Worksheets.Item("LoginData")
$eventSheet = $wb.Worksheets.Item("EventData")

# Order of writes controls shared-string allocation order, so:
# 1) LoginData!A2 -> email address (becomes a hyperlink)
$loginSheet.Range("A2").Value = "raut.ni@northeastern.edu"

# 2) EventData!A3 -> "Code Review"
$eventSheet.Range("A3").Value = "Code Review"

# 3) LoginData!B2 -> "password"
$loginSheet.Range("B2").Value = "password"

# Add the mailto hyperlink on A2, then restore the existing "Hyperlink" cell
# style (Hyperlinks.Add otherwise clones a near-duplicate style record).
$loginSheet.Hyperlinks.Add($loginSheet.Range("A2"), "mailto:raut.ni@northeastern.edu") | Out-Null
$loginSheet.Range("A2").Style = "Hyperlink"

# Update selections per sheet.
$eventSheet.Range("A3").Select() | Out-Null
$loginSheet.Activate() | Out-Null
$loginSheet.Range("B3").Select() | Out-Null
